$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values are not
# reinterpreted/rounded as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '22.108.76'
$ws.Range("E2").Value = '  -0.89%  '

# Row 3
$ws.Range("D3").Value = '1.557.04'
$ws.Range("E3").Value = '  +0.14%  '

# Row 4
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = '1.000'
$ws.Range("E5").Value = '  -0.19%  '

# Row 6
$ws.Range("D6").Value = '292.42'
$ws.Range("E6").Value = '  +1.84%  '

# Row 7
$ws.Range("D7").Value = '0.3955'
$ws.Range("E7").Value = '  +4.80%  '

# Row 8
$ws.Range("D8").Value = '0.3241'
$ws.Range("E8").Value = '  -0.64%  '

# Row 9
$ws.Range("D9").Value = '43.05'
$ws.Range("E9").Value = '  -1.28%  '

# Row 10
$ws.Range("D10").Value = '0.07328'
$ws.Range("E10").Value = '  -0.33%  '

# Row 11
$ws.Range("D11").Value = '1.088'
$ws.Range("E11").Value = '  -3.93%  '

# Row 12
$ws.Range("D12").Value = '1.000'
$ws.Range("E12").Value = '  -0.20%  '

# Row 13
$ws.Range("D13").Value = '19.11'
$ws.Range("E13").Value = '  -5.15%  '

# Row 14
$ws.Range("D14").Value = '5.656'
$ws.Range("E14").Value = '  -2.97%  '

# Row 15
$ws.Range("E15").Value = '  +7.07%  '

# Row 16
$ws.Range("D16").Value = '6.674'
$ws.Range("E16").Value = '  -1.16%  '

# Row 17
$ws.Range("D17").Value = '1.558.87'
$ws.Range("E17").Value = '  -0.29%  '

# Row 18
$ws.Range("D18").Value = '0.06601'
$ws.Range("E18").Value = '  -0.65%  '

# Row 19
$ws.Range("D19").Value = '84.03'
$ws.Range("E19").Value = '  -2.12%  '

# Row 20
$ws.Range("D20").Value = '0.9997'
$ws.Range("E20").Value = '  -0.22%  '

# Row 21
$ws.Range("D21").Value = '6.318'
$ws.Range("E21").Value = '  -0.78%  '

# Row 22
$ws.Range("D22").Value = '15.84'
$ws.Range("E22").Value = '  -1.32%  '

# Row 23
$ws.Range("D23").Value = '11.32'
$ws.Range("E23").Value = '  -2.45%  '

# Row 24
$ws.Range("D24").Value = '22.101.44'
$ws.Range("E24").Value = '  -0.88%  '

# Row 25
$ws.Range("D25").Value = '2.344'
$ws.Range("E25").Value = '  +2.06%  '

# Row 26
$ws.Range("D26").Value = '2.449'
$ws.Range("E26").Value = '  -3.88%  '

# Row 27
$ws.Range("D27").Value = '148.23'
$ws.Range("E27").Value = '  -1.24%  '

# Row 28
$ws.Range("D28").Value = '18.67'
$ws.Range("E28").Value = '  -3.28%  '

# Row 29
$ws.Range("D29").Value = '4.875'
$ws.Range("E29").Value = '  -1.04%  '

# Row 30
$ws.Range("D30").Value = '1.733.91'
$ws.Range("E30").Value = '  -0.48%  '

# Row 31
$ws.Range("D31").Value = '119.31'
$ws.Range("E31").Value = '  -2.40%  '

# Row 32
$ws.Range("D32").Value = '1.033'
$ws.Range("E32").Value = '  -3.93%  '

# Row 33
$ws.Range("D33").Value = '5.714'
$ws.Range("E33").Value = '  -2.94%  '

# Row 34
$ws.Range("D34").Value = '0.08389'
$ws.Range("E34").Value = '  +1.64%  '

# Row 35
$ws.Range("B35").Value = 'WEMIXTOKEN'
$ws.Range("C35").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D35").Value = '1.630'
$ws.Range("E35").Value = '  -14.10%  '

# Row 36
$ws.Range("B36").Value = 'FraxShare'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D36").Value = '9.119'
$ws.Range("E36").Value = '  -2.40%  '

# Row 37
$ws.Range("D37").Value = '0.06167'
$ws.Range("E37").Value = '  -1.89%  '

# Row 38
$ws.Range("D38").Value = '0.02278'
$ws.Range("E38").Value = '  -3.64%  '

# Row 39
$ws.Range("D39").Value = '5.161'
$ws.Range("E39").Value = '  -2.25%  '

# Row 40
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '1.219'
$ws.Range("E40").Value = '  -1.58%  '

# Row 41
$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = '0.2066'
$ws.Range("E41").Value = '  -3.73%  '

# Row 42
$ws.Range("D42").Value = '10.81'
$ws.Range("E42").Value = '  -1.78%  '

# Row 43
$ws.Range("D43").Value = '0.9994'
$ws.Range("E43").Value = '  -0.23%  '

# Row 44
$ws.Range("D44").Value = '0.5868'
$ws.Range("E44").Value = '  -2.88%  '

# Row 45
$ws.Range("D45").Value = '13.18'
$ws.Range("E45").Value = '  -3.60%  '

# Row 46
$ws.Range("D46").Value = '3.768'
$ws.Range("E46").Value = '  +0.78%  '

# Row 47
$ws.Range("D47").Value = '0.5630'
$ws.Range("E47").Value = '  -4.34%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '119.24'
$ws.Range("E48").Value = '  -3.32%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '1.905'
$ws.Range("E49").Value = '  -3.17%  '

# Row 50
$ws.Range("D50").Value = '1.145'

# Row 51
$ws.Range("D51").Value = '0.06866'
$ws.Range("E51").Value = '  -2.91%  '
